# Se procesan de nuevo los datos con las nuevas dimensiones curadas
#
# The dataset's dimension/measure metadata is re-curated:
#  - Several columns that used to be tagged as "iaest-dimension:*" are now
#    correctly tagged as "iaest-measure:*" (and their supporting "dim"/
#    "skos:Concept"/mapping-file rows become "medida"/"xsd:int"/empty).
#  - Column L (municipio-nombre) becomes a proper refArea dimension
#    (sdmx-dimension:refArea / dim / URI-Municipio), mirroring column M
#    (provincia-nombre), which also becomes sdmx-dimension:refArea and
#    keeps its URI-Provincia uri-column marker.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that flip from "dimension" to "measure" metadata.
$measureCols = @("B", "D", "V", "X", "AB")
foreach ($col in $measureCols) {
    $row2 = $ws.Range($col + "2").Value2
    $row2 = $row2 -replace "^iaest-dimension:", "iaest-measure:"
    $ws.Range($col + "2").Value = $row2

    $ws.Range($col + "3").Value = "medida"
    $ws.Range($col + "4").Value = "xsd:int"
    $ws.Range($col + "5").ClearContents()
}

# Column L (municipio-nombre) becomes a refArea dimension, like column M.
$ws.Range("L2").Value = "sdmx-dimension:refArea"
$ws.Range("L3").Value = "dim"
$ws.Range("L4").Value = "URI-Municipio"

# Column O used to be the sdmx-dimension:refPeriod marker; it is now the
# iaest-measure:menos-de-25-anos measure column (already "medida"/"xsd:int"
# style, so only the row2 label and row5 mapping file need to change).
$ws.Range("O2").Value = "iaest-measure:menos-de-25-anos"
$ws.Range("O3").Value = "medida"
$ws.Range("O4").Value = "xsd:int"
$ws.Range("O5").ClearContents()
